# Add average deadline (column J) and deadline std-dev (column K) to each results sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$rows = @(
    @(1, 34.716666666666697, 3.4986276566999801),
    @(2, 43.35, 9.5132824184337892),
    @(3, 47.266666666666701, 4.1451080690320703),
    @(4, 40.133333333333297, 5.8410199558205402),
    @(5, 44.133333333333297, 5.2544304684735001),
    @(6, 52.016666666666701, 5.5218876245223498),
    @(7, 57.6666666666667, 1.4458200027064201),
    @(8, 45.766666666666701, 2.31697250348831),
    @(9, 14.0833333333333, 10.4058192107266),
    @(10, 21.483333333333299, 14.925227950638799),
    @(11, 48.35, 3.1880337705290902),
    @(12, 50.1666666666667, 6.1593711616807596),
    @(13, 43.966666666666697, 7.9040145703106797),
    @(14, 56.966666666666697, 1.54004915909505),
    @(15, 42.45, 2.7085864723072999),
    @(16, 24.05, 10.6809270847539),
    @(17, 33.516666666666701, 15.001120485457401),
    @(18, 51.533333333333303, 5.8843851590939398),
    @(19, 32.283333333333303, 6.0983233779948796),
    @(20, 52.466666666666697, 5.7709578756268796),
    @(21, 48.066666666666698, 3.3133866471401001),
    @(22, 53.616666666666703, 5.5206597083880498),
    @(23, 25.533333333333299, 13.853307341345801),
    @(24, 51.783333333333303, 2.4570318780425899),
    @(25, 2.6666666666666701, 0.47538268854152799),
    @(26, 41.133333333333297, 6.02105533479246),
    @(27, 56.566666666666698, 1.3575484717340101),
    @(28, 41.383333333333297, 2.3512828569626301)
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $ws.Cells.Item($r, 11).Value = $row[2]
}
$ws.Range("J1:K28").Select()

$ws = $wb.Worksheets.Item(2)
$rows = @(
    @(1, 99.383333333333297, 32.360831978298798),
    @(2, 94.8, 17.4810308814431),
    @(3, 69.75, 46.114310804965399),
    @(4, 60.05, 18.742162203646501),
    @(5, 95.266666666666694, 16.491257742237401),
    @(6, 89.366666666666703, 10.697642171179499),
    @(7, 81.866666666666703, 12.512321610732201),
    @(8, 115.433333333333, 28.9835527098321),
    @(9, 67.4166666666667, 17.242332284202998),
    @(10, 103.26666666666701, 29.105807368957699),
    @(11, 36.200000000000003, 34.679549974598501),
    @(12, 66.900000000000006, 34.742003103750598),
    @(13, 62.1666666666667, 24.257912681919201),
    @(14, 73.099999999999994, 18.417659175137501),
    @(15, 93.983333333333306, 49.518184184920898),
    @(16, 41.8333333333333, 11.1601085388481),
    @(17, 108.566666666667, 16.031078573383301),
    @(18, 56.6, 49.901326363949103),
    @(19, 99.8, 14.072210865534499),
    @(20, 97.266666666666694, 18.121872416696),
    @(21, 136.38333333333301, 18.4823247369917),
    @(22, 72.766666666666694, 30.360750011205901),
    @(23, 105.383333333333, 29.4498088777814),
    @(24, 109.633333333333, 45.5113942397914),
    @(25, 87.9, 12.942768414728301),
    @(26, 169.833333333333, 2.9754171715295299),
    @(27, 67.216666666666697, 19.692631046006099),
    @(28, 108.2, 14.616916134675501),
    @(29, 99.9166666666667, 28.9279008282261),
    @(30, 118.65, 27.982939717875201),
    @(31, 75.7, 15.714939453714001),
    @(32, 133.51666666666699, 30.437761820356801),
    @(33, 132.01666666666699, 47.110577338588598),
    @(34, 98.1666666666667, 45.685426473771599),
    @(35, 97.366666666666703, 32.678177379348398),
    @(36, 41.0833333333333, 40.079705193047303)
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $ws.Cells.Item($r, 11).Value = $row[2]
}
$ws.Range("J1:K36").Select()

$ws = $wb.Worksheets.Item(3)
$rows = @(
    @(1, 99.616666666666703, 18.8185732871802),
    @(2, 119.916666666667, 25.333231117843098),
    @(3, 102.783333333333, 8.5788302048482006),
    @(4, 128.53333333333299, 7.8620881054927798),
    @(5, 104.933333333333, 10.881280250366),
    @(6, 106.35, 8.35205242207069),
    @(7, 162.25, 5.4169382808910704),
    @(8, 101.4, 7.4017863166918101),
    @(9, 143.80000000000001, 24.130857943403999),
    @(10, 97.183333333333294, 37.0054927973557),
    @(11, 129.19999999999999, 10.110910368239701),
    @(12, 104.51666666666701, 13.011717882735001),
    @(13, 101.666666666667, 7.4188643350986601),
    @(14, 165.1, 4.5349864051521998),
    @(15, 79.75, 14.240816249696801),
    @(16, 115.683333333333, 18.450567642424499),
    @(17, 128.38333333333301, 7.0086911340441702),
    @(18, 127.216666666667, 7.7088803440053599),
    @(19, 129.01666666666699, 9.2415635676621992),
    @(20, 163.05000000000001, 3.4268382423136101),
    @(21, 125.166666666667, 8.7936160716635303),
    @(22, 126.916666666667, 31.6523979152498),
    @(23, 104.133333333333, 9.05438683202539),
    @(24, 161.23333333333301, 6.4686117407916299),
    @(25, 102.116666666667, 11.8236930573372),
    @(26, 113.166666666667, 12.802497990714899),
    @(27, 110.583333333333, 12.110076581188901),
    @(28, 162.333333333333, 3.7356127020727201),
    @(29, 100.883333333333, 6.1205174996299796),
    @(30, 148.6, 32.150177267086498),
    @(31, 147.36666666666699, 31.586827691552401),
    @(32, 162.38333333333301, 4.7159723790779404),
    @(33, 120.883333333333, 20.082956486606701),
    @(34, 102.716666666667, 33.973614720473002)
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $ws.Cells.Item($r, 11).Value = $row[2]
}
$ws.Range("J1:K34").Select()

$ws = $wb.Worksheets.Item(4)
$rows = @(
    @(1, 158.19999999999999, 23.895819079737102),
    @(2, 135.816666666667, 23.965718218269899),
    @(3, 118.2, 8.9457877598380193),
    @(4, 141.53333333333299, 6.2394770170526099),
    @(5, 127.166666666667, 16.583805319348599),
    @(6, 121.933333333333, 6.2649594982030301),
    @(7, 165.98333333333301, 2.92558362519731),
    @(8, 119.51666666666701, 5.6583271968547999),
    @(9, 164.75, 20.380802698717901),
    @(10, 104.2, 30.643272826046601),
    @(11, 146.51666666666699, 8.9546724223776302),
    @(12, 127.116666666667, 28.627408738365499),
    @(13, 118.083333333333, 7.8832404648081598),
    @(14, 169.833333333333, 3.87589974974283),
    @(15, 108, 11.6123838046767),
    @(16, 119.25, 24.9695577365584),
    @(17, 150.71666666666701, 15.7792952374968),
    @(18, 145.183333333333, 8.8958029842428097),
    @(19, 138.85, 6.5037798657682799),
    @(20, 166.25, 2.6013360322231298),
    @(21, 139.36666666666699, 6.9647214325318396),
    @(22, 132.01666666666699, 33.8393471040648),
    @(23, 129.26666666666699, 32.558092366277798),
    @(24, 125.833333333333, 16.802105060476102),
    @(25, 168.86666666666699, 2.9998116701715101),
    @(26, 110.9, 11.563677791593699),
    @(27, 115.283333333333, 16.0434702559361),
    @(28, 113.3, 15.8042045467633),
    @(29, 166.083333333333, 3.0437023990175698),
    @(30, 115.383333333333, 5.5298624332865201),
    @(31, 149.683333333333, 32.886807239309199),
    @(32, 148.94999999999999, 32.366767325740803),
    @(33, 166.38333333333301, 2.7927974352115701),
    @(34, 139.03333333333299, 18.022553353466101),
    @(35, 119.116666666667, 24.030412698563399)
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $ws.Cells.Item($r, 11).Value = $row[2]
}
$ws.Range("J1:K35").Select()
